$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F11").Value = 21
$ws1.Range("F18").Value = 125
$ws1.Range("F22").Value = 6769
$ws1.Range("F23").Value = 235
$ws1.Range("F26").Value = 1258
$ws1.Range("F35").Value = 84
$ws1.Range("F37").Value = 4688
$ws1.Range("F46").Value = 1007

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F17").Value = 125
$ws4.Range("F21").Value = 6769
$ws4.Range("F22").Value = 235
$ws4.Range("F25").Value = 1258
$ws4.Range("F36").Value = 84
$ws4.Range("F38").Value = 4688
$ws4.Range("F47").Value = 1007
